# Apply "Add data for 2025-07-08" update to cta-violent-crime-ytd.xlsx
# This applies a batch of retroactive/incremental count corrections across
# the Citywide Totals sheet, the By Neighborhood summary sheet, and several
# individual per-neighborhood sheets, plus inserts one brand-new
# crime-category row (Aggravated Battery) into the Morgan Park sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Citywide Totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("D2").Value = 49
$ws.Range("J2").Value = 67
$ws.Range("K2").Value = 74
$ws.Range("D3").Value = 74
$ws.Range("F3").Value = 72
$ws.Range("K3").Value = 121
$ws.Range("L3").Value = 128
$ws.Range("K4").Value = 15
$ws.Range("C6").Value = 243
$ws.Range("E6").Value = 225
$ws.Range("F6").Value = 272
$ws.Range("G6").Value = 252
$ws.Range("J6").Value = 217
$ws.Range("K6").Value = 254
$ws.Range("L6").Value = 266
$ws.Range("C7").Value = 329
$ws.Range("D7").Value = 357
$ws.Range("E7").Value = 343
$ws.Range("F7").Value = 391
$ws.Range("G7").Value = 367
$ws.Range("J7").Value = 399
$ws.Range("K7").Value = 467
$ws.Range("L7").Value = 484

# ---------------------------------------------------------------
# Garfield Park
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("E6").Value = 18
$ws.Range("E7").Value = 26

# ---------------------------------------------------------------
# Grand Crossing
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 35

# ---------------------------------------------------------------
# Englewood
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 19
$ws.Range("L5").Value = 26
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 26
$ws.Range("L6").Value = 41

# ---------------------------------------------------------------
# By Neighborhood (summary/pivot sheet)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 5
$ws.Range("C7").Value = 25
$ws.Range("F7").Value = 28
$ws.Range("J7").Value = 27
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 26
$ws.Range("L26").Value = 41
$ws.Range("K28").Value = 9
$ws.Range("E30").Value = 26
$ws.Range("L34").Value = 35
$ws.Range("K40").Value = 2
$ws.Range("C51").Value = 28
$ws.Range("D51").Value = 49
$ws.Range("F51").Value = 39
$ws.Range("G51").Value = 46
$ws.Range("L51").Value = 55
$ws.Range("F52").Value = 3
$ws.Range("D57").Value = 3
$ws.Range("J63").Value = 3
$ws.Range("G70").Value = 4
$ws.Range("K70").Value = 5
$ws.Range("K72").Value = 10
$ws.Range("C96").Value = 329
$ws.Range("D96").Value = 357
$ws.Range("E96").Value = 343
$ws.Range("F96").Value = 391
$ws.Range("G96").Value = 367
$ws.Range("J96").Value = 399
$ws.Range("K96").Value = 467
$ws.Range("L96").Value = 484

# ---------------------------------------------------------------
# Loop
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D2").Value = 8
$ws.Range("C6").Value = 17
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = 31
$ws.Range("L6").Value = 25
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = 49
$ws.Range("F7").Value = 39
$ws.Range("G7").Value = 46
$ws.Range("L7").Value = 55

# ---------------------------------------------------------------
# North Lawndale (new Aggravated Assault value for 2023)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 1
$ws.Range("J6").Value = 3

# ---------------------------------------------------------------
# River North (new Aggravated Assault value for 2024)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 1
$ws.Range("K6").Value = 10

# ---------------------------------------------------------------
# Printers Row
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("G4").Value = 3
$ws.Range("J4").Value = 5
$ws.Range("G5").Value = 4
$ws.Range("J5").Value = 5

# ---------------------------------------------------------------
# Hyde Park (new Aggravated Battery value for 2024)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I3").Value = 1
$ws.Range("I5").Value = 2

# ---------------------------------------------------------------
# Auburn Gresham (new Aggravated Battery value for 2024)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 1
$ws.Range("K6").Value = 5

# ---------------------------------------------------------------
# Albany Park (new Criminal Sexual Assault value for 2024)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("H4").Value = 1
$ws.Range("H6").Value = 4

# ---------------------------------------------------------------
# Gage Park
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 5
$ws.Range("J5").Value = 9

# ---------------------------------------------------------------
# Lower West Side (new Aggravated Battery value for 2019)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("F3").Value = 1
$ws.Range("F5").Value = 3

# ---------------------------------------------------------------
# Austin
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("C5").Value = 17
$ws.Range("F5").Value = 19
$ws.Range("J5").Value = 13
$ws.Range("C6").Value = 25
$ws.Range("F6").Value = 28
$ws.Range("J6").Value = 27

# ---------------------------------------------------------------
# Morgan Park - a brand-new "Aggravated Battery" crime_category row
# needs to be inserted (in crime_category/shared-string order) between
# the existing "Aggravated Assault" row and the existing "Robbery" row,
# pushing "Robbery" and "Total" down by one row, and updating the
# "Total" row's 2017 figure to include the new record.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Rows.Item(3).Insert()

# Copy the row-label formatting (bold/border/center) from an existing
# label cell onto the newly inserted label cell.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122)

$ws.Cells.Item(3, 1).Value = "Aggravated Battery"
$ws.Cells.Item(3, 3).Value = 1

# Update the (now shifted-down) Total row's 2017 count
$ws.Range("C5").Value = 3
